# The workbook's weekly "Hortaliza, Femacal de La Calera - Alcachofa" data set
# gained one additional daily record. A new row is inserted at row 120
# (shifting the existing rows 120-176 down to 121-177) and populated with
# the new observation's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 120, pushing rows 120-176 down to 121-177.
$ws.Rows("120:120").Insert()

# Populate the newly inserted row 120 with the new record.
$ws.Cells.Item(120, 1).Value = 3
$ws.Cells.Item(120, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(120, 3).Value = "Coquimbo"
$ws.Cells.Item(120, 4).Value = 44455
$ws.Cells.Item(120, 5).Value = 5
$ws.Cells.Item(120, 6).Value = 100112013
$ws.Cells.Item(120, 7).Value = "Alcachofa"
$ws.Cells.Item(120, 8).Value = "Española"
$ws.Cells.Item(120, 9).Value = "Extra"
$ws.Cells.Item(120, 10).Value = 115
$ws.Cells.Item(120, 11).Value = 11500
$ws.Cells.Item(120, 12).Value = 12000
$ws.Cells.Item(120, 13).Value = 11739
$ws.Cells.Item(120, 14).Value = "`$/caja 30 unidades"
$ws.Cells.Item(120, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(120, 16).Value = 391
$ws.Cells.Item(120, 17).Value = 30
$ws.Cells.Item(120, 18).Value = "Hortaliza"
